# Colony 04 revision initial complete, overworld 2 revision topography maps added
# Insert a new "Swamp" topography entry above the existing "warning" row (row 39)
# on the language sheet, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 39 (currently "warning" / "WARNING"),
# pushing it and everything below down by one row.
$ws.Rows.Item(39).Insert()

# Populate the new row with the Swamp topography key/value pair.
$ws.Cells.Item(39, 1).Value = "topography_swamp"
$ws.Cells.Item(39, 2).Value = "Swamp"

# Match the wrap-text formatting used by the other Value cells in column B.
$ws.Cells.Item(39, 2).WrapText = $true

# Update the view so the active selection reflects where the edit took place
# (the row that used to be "warning" is now row 40).
$ws.Activate()
$ws.Range("B40").Select()
